$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1900.0312
$ws.Range("I15").Value = 1900.0312
$ws.Range("K15").Value = 5700.0936
$ws.Range("M15").Value = -5531.0936
$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 300
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H43").Value = 2865
$ws.Range("J43").Value = 3650
$ws.Range("L43").Value = 3650
$ws.Range("N43").Value = -3788
$ws.Range("H69").Value = 33402790
$ws.Range("I69").Value = 7833
$ws.Range("J69").Value = 41751532
$ws.Range("K69").Value = 23499
$ws.Range("L69").Value = 125254596
$ws.Range("M69").Value = -22625
$ws.Range("N69").Value = -125256344
$ws.Range("H72").Value = 33402790
$ws.Range("I72").Value = 7833
$ws.Range("J72").Value = 41751532
$ws.Range("K72").Value = 70497
$ws.Range("L72").Value = 375763788
$ws.Range("M72").Value = -66129
$ws.Range("N72").Value = -375772524
$ws.Range("H88").Value = 4440.4165
$ws.Range("I88").Value = 5956
$ws.Range("J88").Value = 1409.25
$ws.Range("K88").Value = 5956
$ws.Range("L88").Value = 1409.25
$ws.Range("M88").Value = -5550
$ws.Range("N88").Value = -2221.25
$ws.Range("H91").Value = 4440.4165
$ws.Range("I91").Value = 5956
$ws.Range("J91").Value = 1409.25
$ws.Range("K91").Value = 5956
$ws.Range("L91").Value = 1409.25
$ws.Range("M91").Value = -4552
$ws.Range("N91").Value = -4217.25
$ws.Range("H92").Value = 8334556.5
$ws.Range("I92").Value = 1127.8889
$ws.Range("K92").Value = 1127.8889
$ws.Range("M92").Value = 120.1111000000001
$ws.Range("H98").Value = 1689.909
$ws.Range("J98").Value = 895
$ws.Range("L98").Value = 895
$ws.Range("N98").Value = -3891
$ws.Range("H106").Value = 2322.5557
$ws.Range("I106").Value = 2050.375
$ws.Range("K106").Value = 2050.375
$ws.Range("M106").Value = -1419.375
$ws.Range("H122").Value = 1689.909
$ws.Range("J122").Value = 895
$ws.Range("L122").Value = 2685
$ws.Range("N122").Value = -7585
$ws.Range("H123").Value = 59999.91
$ws.Range("J123").Value = 59999.91
$ws.Range("L123").Value = 59999.91
$ws.Range("N123").Value = -69799.91
$ws.Range("H125").Value = 13086.125
$ws.Range("I125").Value = 664.7143
$ws.Range("J125").Value = 100036
$ws.Range("K125").Value = 5982.428699999999
$ws.Range("L125").Value = 900324
$ws.Range("M125").Value = -3522.428699999999
$ws.Range("N125").Value = -905244
$ws.Range("H132").Value = 56749.445
$ws.Range("I132").Value = 1354.1428
$ws.Range("J132").Value = 250633
$ws.Range("K132").Value = 4062.4284
$ws.Range("L132").Value = 751899
$ws.Range("M132").Value = -1532.4284
$ws.Range("N132").Value = -756959
$ws.Range("H135").Value = 7821
$ws.Range("I135").Value = 2541.5334
$ws.Range("J135").Value = 17720
$ws.Range("K135").Value = 22873.8006
$ws.Range("L135").Value = 159480
$ws.Range("M135").Value = -20338.8006
$ws.Range("N135").Value = -164550
$ws.Range("H137").Value = 3046.8718
$ws.Range("I137").Value = 2736.762
$ws.Range("J137").Value = 4349.3335
$ws.Range("K137").Value = 8210.286
$ws.Range("L137").Value = 13048.0005
$ws.Range("M137").Value = -5660.286
$ws.Range("N137").Value = -18148.0005
$ws.Range("H138").Value = 1900.3651
$ws.Range("I138").Value = 747.45
$ws.Range("J138").Value = 3905.4348
$ws.Range("K138").Value = 2242.35
$ws.Range("L138").Value = 11716.3044
$ws.Range("M138").Value = 2897.65
$ws.Range("N138").Value = -21996.3044

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3092.389
$ws.Range("I45").Value = 3095.1538
$ws.Range("K45").Value = 3095.1538
$ws.Range("M45").Value = -2718.1538
$ws.Range("H61").Value = 4854.8667
$ws.Range("I61").Value = 3594.6
$ws.Range("K61").Value = 3594.6
$ws.Range("M61").Value = -3382.6
$ws.Range("H74").Value = 1816.7084
$ws.Range("I74").Value = 1483.85
$ws.Range("J74").Value = 3481
$ws.Range("K74").Value = 1483.85
$ws.Range("L74").Value = 3481
$ws.Range("M74").Value = -609.8499999999999
$ws.Range("N74").Value = -5229
$ws.Range("H76").Value = 66713.8
$ws.Range("J76").Value = 66713.8
$ws.Range("L76").Value = 66713.8
$ws.Range("N76").Value = -67389.8
$ws.Range("H77").Value = 1816.7084
$ws.Range("I77").Value = 1483.85
$ws.Range("J77").Value = 3481
$ws.Range("K77").Value = 7419.25
$ws.Range("L77").Value = 17405
$ws.Range("M77").Value = -3051.25
$ws.Range("N77").Value = -26141
$ws.Range("H79").Value = 66713.8
$ws.Range("J79").Value = 66713.8
$ws.Range("L79").Value = 66713.8
$ws.Range("N79").Value = -69053.8
$ws.Range("H110").Value = 1228.7059
$ws.Range("I110").Value = 1259.2
$ws.Range("K110").Value = 1259.2
$ws.Range("M110").Value = 785.8
$ws.Range("H122").Value = 2395.7727
$ws.Range("I122").Value = 2381.75
$ws.Range("J122").Value = 2433.1667
$ws.Range("K122").Value = 7145.25
$ws.Range("L122").Value = 7299.500100000001
$ws.Range("M122").Value = -4695.25
$ws.Range("N122").Value = -12199.5001
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1910.2142
$ws.Range("I132").Value = 1888.3704
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5665.1112
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3135.1112
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 4854.8667
$ws.Range("I136").Value = 3594.6
$ws.Range("K136").Value = 10783.8
$ws.Range("M136").Value = -8233.799999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 892.3200000000001
$ws.Range("I20").Value = 875.9375
$ws.Range("J20").Value = 921.44446
$ws.Range("K20").Value = 875.9375
$ws.Range("L20").Value = 921.44446
$ws.Range("M20").Value = -628.9375
$ws.Range("N20").Value = -1415.44446
$ws.Range("H105").Value = 4039.7222
$ws.Range("I105").Value = 2027.3334
$ws.Range("J105").Value = 5045.9165
$ws.Range("K105").Value = 2027.3334
$ws.Range("L105").Value = 5045.9165
$ws.Range("M105").Value = -280.3334
$ws.Range("N105").Value = -8539.916499999999
$ws.Range("H134").Value = 2571.7407
$ws.Range("I134").Value = 2202.162
$ws.Range("J134").Value = 3376.1177
$ws.Range("K134").Value = 6606.485999999999
$ws.Range("L134").Value = 10128.3531
$ws.Range("M134").Value = -4071.485999999999
$ws.Range("N134").Value = -15198.3531
$ws.Range("H141").Value = 200780
$ws.Range("J141").Value = 200780
$ws.Range("L141").Value = 200780
$ws.Range("N141").Value = -211140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 257.66666
$ws.Range("I7").Value = 182.42857
$ws.Range("J7").Value = 323.5
$ws.Range("K7").Value = 182.42857
$ws.Range("L7").Value = 323.5
$ws.Range("M7").Value = -69.42857000000001
$ws.Range("N7").Value = -549.5
$ws.Range("H16").Value = 3430.8572
$ws.Range("I16").Value = 3776.5908
$ws.Range("K16").Value = 3776.5908
$ws.Range("M16").Value = -3489.5908
$ws.Range("H31").Value = 2163.5
$ws.Range("I31").Value = 1991.8572
$ws.Range("J31").Value = 2884.4
$ws.Range("K31").Value = 1991.8572
$ws.Range("L31").Value = 2884.4
$ws.Range("M31").Value = -1696.8572
$ws.Range("N31").Value = -3474.4
$ws.Range("H34").Value = 2163.5
$ws.Range("I34").Value = 1991.8572
$ws.Range("J34").Value = 2884.4
$ws.Range("K34").Value = 1991.8572
$ws.Range("L34").Value = 2884.4
$ws.Range("M34").Value = -1789.8572
$ws.Range("N34").Value = -3288.4
$ws.Range("H58").Value = 2688
$ws.Range("I58").Value = 1475.091
$ws.Range("J58").Value = 4092.4211
$ws.Range("K58").Value = 1475.091
$ws.Range("L58").Value = 4092.4211
$ws.Range("M58").Value = -1272.091
$ws.Range("N58").Value = -4498.4211
$ws.Range("H88").Value = 29771.25
$ws.Range("J88").Value = 29771.25
$ws.Range("L88").Value = 29771.25
$ws.Range("N88").Value = -30583.25
$ws.Range("H91").Value = 29771.25
$ws.Range("J91").Value = 29771.25
$ws.Range("L91").Value = 29771.25
$ws.Range("N91").Value = -32579.25
$ws.Range("H99").Value = 17265644
$ws.Range("I99").Value = 2442033.5
$ws.Range("K99").Value = 2442033.5
$ws.Range("M99").Value = -2440535.5
$ws.Range("H105").Value = 17998.166
$ws.Range("I105").Value = 1500
$ws.Range("J105").Value = 26247.25
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 26247.25
$ws.Range("M105").Value = 247
$ws.Range("N105").Value = -29741.25
$ws.Range("H107").Value = 8285.6
$ws.Range("I107").Value = 1562.25
$ws.Range("J107").Value = 15969.429
$ws.Range("K107").Value = 1562.25
$ws.Range("L107").Value = 15969.429
$ws.Range("M107").Value = 357.75
$ws.Range("N107").Value = -19809.429
$ws.Range("H113").Value = 3430.8572
$ws.Range("I113").Value = 3776.5908
$ws.Range("K113").Value = 3776.5908
$ws.Range("M113").Value = -1606.5908
$ws.Range("H121").Value = 19999
$ws.Range("J121").Value = 19999
$ws.Range("L121").Value = 19999
$ws.Range("N121").Value = -22619
$ws.Range("H126").Value = 17265644
$ws.Range("I126").Value = 2442033.5
$ws.Range("K126").Value = 7326100.5
$ws.Range("M126").Value = -7323630.5
$ws.Range("H132").Value = 900.1
$ws.Range("I132").Value = 900.1
$ws.Range("K132").Value = 2700.3
$ws.Range("M132").Value = -170.3000000000002
$ws.Range("H134").Value = 2498.8958
$ws.Range("I134").Value = 2057.186
$ws.Range("K134").Value = 6171.558000000001
$ws.Range("M134").Value = -3636.558000000001
$ws.Range("H136").Value = 2688
$ws.Range("I136").Value = 1475.091
$ws.Range("J136").Value = 4092.4211
$ws.Range("K136").Value = 4425.272999999999
$ws.Range("L136").Value = 12277.2633
$ws.Range("M136").Value = -1875.272999999999
$ws.Range("N136").Value = -17377.2633

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27189748
$ws.Range("I4").Value = 33533826
$ws.Range("J4").Value = 843.8570999999999
$ws.Range("K4").Value = 100601478
$ws.Range("L4").Value = 2531.5713
$ws.Range("M4").Value = -100601366
$ws.Range("N4").Value = -2755.5713
$ws.Range("H39").Value = 4081.2666
$ws.Range("J39").Value = 4081.2666
$ws.Range("L39").Value = 12243.7998
$ws.Range("N39").Value = -12831.7998
$ws.Range("H55").Value = 1252788
$ws.Range("J55").Value = 1431714.9
$ws.Range("L55").Value = 4295144.699999999
$ws.Range("N55").Value = -4295498.699999999
$ws.Range("H107").Value = 561.375
$ws.Range("J107").Value = 558.7143
$ws.Range("L107").Value = 1676.1429
$ws.Range("N107").Value = -5516.1429
$ws.Range("H122").Value = 675.3684
$ws.Range("I122").Value = 709.0714
$ws.Range("J122").Value = 655.7083
$ws.Range("K122").Value = 6381.6426
$ws.Range("L122").Value = 5901.3747
$ws.Range("M122").Value = -3931.6426
$ws.Range("N122").Value = -10801.3747
$ws.Range("H131").Value = 2568.2083
$ws.Range("I131").Value = 1284.1875
$ws.Range("J131").Value = 5136.25
$ws.Range("K131").Value = 3852.5625
$ws.Range("L131").Value = 15408.75
$ws.Range("M131").Value = 1187.4375
$ws.Range("N131").Value = -25488.75
$ws.Range("H132").Value = 1009.6
$ws.Range("I132").Value = 1074.5
$ws.Range("J132").Value = 750
$ws.Range("K132").Value = 9670.5
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -7140.5
$ws.Range("N132").Value = -11810
$ws.Range("H139").Value = 6873.9
$ws.Range("I139").Value = 5147.8
$ws.Range("K139").Value = 15443.4
$ws.Range("M139").Value = -10303.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4000241.5
$ws.Range("I2").Value = 222.5
$ws.Range("J2").Value = 5263405.5
$ws.Range("K2").Value = 222.5
$ws.Range("L2").Value = 5263405.5
$ws.Range("M2").Value = -109.5
$ws.Range("N2").Value = -5263631.5
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -4827
$ws.Range("N21").ClearContents()
$ws.Range("H24").Value = 8765396
$ws.Range("I24").Value = 17500000
$ws.Range("J24").Value = 30791.375
$ws.Range("K24").Value = 17500000
$ws.Range("L24").Value = 30791.375
$ws.Range("M24").Value = -17499827
$ws.Range("N24").Value = -31137.375
$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 5000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -4895
$ws.Range("N30").ClearContents()
$ws.Range("H97").Value = 12500503
$ws.Range("J97").Value = 25000408
$ws.Range("L97").Value = 25000408
$ws.Range("N97").Value = -25001400
$ws.Range("H102").Value = 1637.409
$ws.Range("I102").Value = 1017.37836
$ws.Range("J102").Value = 4914.7144
$ws.Range("K102").Value = 1017.37836
$ws.Range("L102").Value = 4914.7144
$ws.Range("M102").Value = 604.62164
$ws.Range("N102").Value = -8158.7144
$ws.Range("H107").Value = 1139.3636
$ws.Range("I107").Value = 1259.2222
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1259.2222
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 660.7778000000001
$ws.Range("N107").Value = -4440
$ws.Range("H113").Value = 3646.5
$ws.Range("I113").Value = 1775.9445
$ws.Range("J113").Value = 7013.5
$ws.Range("K113").Value = 1775.9445
$ws.Range("L113").Value = 7013.5
$ws.Range("M113").Value = 394.0554999999999
$ws.Range("N113").Value = -11353.5
$ws.Range("H123").Value = 42000
$ws.Range("J123").Value = 42000
$ws.Range("L123").Value = 42000
$ws.Range("N123").Value = -46900
$ws.Range("H132").Value = 2386.6938
$ws.Range("I132").Value = 2421.6924
$ws.Range("J132").Value = 2250.2
$ws.Range("K132").Value = 7265.0772
$ws.Range("L132").Value = 6750.599999999999
$ws.Range("M132").Value = -4735.0772
$ws.Range("N132").Value = -11810.6
$ws.Range("H134").Value = 107500
$ws.Range("J134").Value = 107500
$ws.Range("L134").Value = 322500
$ws.Range("N134").Value = -327570
$ws.Range("H136").Value = 18294.555
$ws.Range("J136").Value = 18294.555
$ws.Range("L136").Value = 54883.665
$ws.Range("N136").Value = -59983.665

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10873
$ws.Range("I40").Value = 14340.25
$ws.Range("K40").Value = 14340.25
$ws.Range("M40").Value = -14204.25
$ws.Range("H61").Value = 2178.1428
$ws.Range("I61").Value = 2611.75
$ws.Range("J61").Value = 1600
$ws.Range("K61").Value = 2611.75
$ws.Range("L61").Value = 1600
$ws.Range("M61").Value = -2409.75
$ws.Range("N61").Value = -2004
$ws.Range("H93").Value = 14708930
$ws.Range("I93").Value = 3504.25
$ws.Range("J93").Value = 50001950
$ws.Range("K93").Value = 3504.25
$ws.Range("L93").Value = 50001950
$ws.Range("M93").Value = -2256.25
$ws.Range("N93").Value = -50004446
$ws.Range("H100").Value = 66158.836
$ws.Range("I100").Value = 279739.75
$ws.Range("K100").Value = 279739.75
$ws.Range("M100").Value = -279198.75
$ws.Range("H113").Value = 2178.1428
$ws.Range("I113").Value = 2611.75
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 2611.75
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = -441.75
$ws.Range("N113").Value = -5940
$ws.Range("H132").Value = 5954.8066
$ws.Range("I132").Value = 4797.6
$ws.Range("K132").Value = 14392.8
$ws.Range("M132").Value = -11862.8
$ws.Range("H136").Value = 3423.5557
$ws.Range("I136").Value = 2776.5881
$ws.Range("J136").Value = 4523.4
$ws.Range("K136").Value = 8329.764299999999
$ws.Range("L136").Value = 13570.2
$ws.Range("M136").Value = -5779.764299999999
$ws.Range("N136").Value = -18670.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 398
$ws.Range("I7").Value = 447.5
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 447.5
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -334.5
$ws.Range("N7").Value = -426
$ws.Range("H56").Value = 9999
$ws.Range("I56").Value = 9999
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 9999
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -9285
$ws.Range("N56").ClearContents()
$ws.Range("H81").Value = 13893539
$ws.Range("I81").Value = 1066.3334
$ws.Range("K81").Value = 2132.6668
$ws.Range("M81").Value = -1071.6668
$ws.Range("H84").Value = 13893539
$ws.Range("I84").Value = 1066.3334
$ws.Range("K84").Value = 10663.334
$ws.Range("M84").Value = -5359.333999999999
$ws.Range("H107").Value = 630.2083
$ws.Range("I107").Value = 706.8889
$ws.Range("K107").Value = 2120.6667
$ws.Range("M107").Value = -200.6667000000002
$ws.Range("H122").Value = 283246.56
$ws.Range("I122").Value = 2130.8708
$ws.Range("K122").Value = 6392.6124
$ws.Range("M122").Value = -3942.6124
$ws.Range("H132").Value = 1181.8928
$ws.Range("I132").Value = 969.42554
$ws.Range("J132").Value = 2291.4443
$ws.Range("K132").Value = 2908.27662
$ws.Range("L132").Value = 6874.3329
$ws.Range("M132").Value = -378.2766199999996
$ws.Range("N132").Value = -11934.3329
$ws.Range("H136").Value = 1332.7
$ws.Range("I136").Value = 854.62964
$ws.Range("J136").Value = 2325.6155
$ws.Range("K136").Value = 2563.88892
$ws.Range("L136").Value = 6976.8465
$ws.Range("M136").Value = -13.88891999999987
$ws.Range("N136").Value = -12076.8465
